$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scores")

$ws.Range("C2").Value = "Very Glib"
$ws.Range("C3").Value = "Very Grandiose"
$ws.Range("C4").Value = "Very Conniving"
$ws.Range("C5").Value = "Very Deceptive"
$ws.Range("C6").Value = "Very Unremorseful"
$ws.Range("C7").Value = "Very Callous"
$ws.Range("C8").Value = "Very Inexpressive"
$ws.Range("C9").Value = "Very Irresponsible"
$ws.Range("C10").Value = "Very Sensation Seeking"
$ws.Range("C11").Value = "Very Unrealistic"
$ws.Range("C12").Value = "Very Impulsive"
$ws.Range("C13").Value = "Very Irresponsible"
$ws.Range("C14").Value = "Very Parasitic"
$ws.Range("C15").Value = "Very Noncommittal"
$ws.Range("C16").Value = "Very Promiscuous"
$ws.Range("C17").Value = "Very Emotionally Controlled"
$ws.Range("C18").Value = "Very Problematic"
$ws.Range("C19").Value = "Very Delinquent"
$ws.Range("C20").Value = "Very Noncompliant"
$ws.Range("C21").Value = "Very Versatile"
